$d = $word.ActiveDocument

# --- Change 1: clear the "Post condiciones" summary-table row ------------
# Row 8 of table 2 holds the hyperlinked text "La postulación queda
# registrada en el sistema." in column 1 and the page number "5" in
# column 2. Both runs (and the hyperlink wrapping the first one) are
# removed, leaving the paragraphs in place but empty.
$tbl = $d.Tables(2)

$cell1 = $tbl.Cell(8, 1)
$r1 = $d.Range($cell1.Range.Start, $cell1.Range.End)
$r1.Delete()

$cell2 = $tbl.Cell(8, 2)
$r2 = $d.Range($cell2.Range.Start, $cell2.Range.End)
$r2.Delete()

# --- Change 2: add a new precondition bullet -------------------------------
# After the "... está logeado a su cuenta." paragraph, insert a new list
# paragraph (same list formatting) with the additional precondition text.
$rngFind = $d.Content
$found = $rngFind.Find.Execute("está logeado a su cuenta.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rngFind.Paragraphs(1)
$para.Range.InsertParagraphAfter()
$nextPara = $para.Next()
$nextPara.Range.InsertBefore("El usuario tiene los datos que va a querer actualizar o modificar.")
